$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.320.95"
$ws.Range("E2").Value = "  +2.40%  "

# Row 3
$ws.Range("D3").Value = "3.390.09"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.84"
$ws.Range("E5").Value = "  +1.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.70"
$ws.Range("E6").Value = "  +3.00%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  +1.45%  "

# Row 9
$ws.Range("E9").Value = "  +7.94%  "

# Row 10
$ws.Range("E10").Value = "  +2.59%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.68"
$ws.Range("E11").Value = "  +4.08%  "

# Row 12
$ws.Range("E12").Value = "  +3.90%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "680.19"
$ws.Range("E13").Value = "  -1.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.67"
$ws.Range("E14").Value = "  +3.06%  "

# Row 15
$ws.Range("D15").Value = "3.931.45"
$ws.Range("E15").Value = "  +1.40%  "

# Row 16
$ws.Range("D16").Value = "69.381.62"
$ws.Range("E16").Value = "  +2.40%  "

# Row 17
$ws.Range("D17").Value = "3.398.42"

# Row 18
$ws.Range("E18").Value = "  +1.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.75"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.32"
$ws.Range("E20").Value = "  +2.61%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.906"
$ws.Range("E21").Value = "  +1.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.42"
$ws.Range("E22").Value = "  -0.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.17"
$ws.Range("E23").Value = "  +1.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.12"
$ws.Range("E24").Value = "  +1.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").Value = "  +0.61%  "

# Row 26
$ws.Range("E26").Value = "  +2.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  +1.86%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.95"
$ws.Range("E28").Value = "  +3.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.76"
$ws.Range("E29").Value = "  +2.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -0.87%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.16"
$ws.Range("E31").Value = "  +1.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.65"
$ws.Range("E32").Value = "  +11.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "556.08"
$ws.Range("E33").Value = "  -2.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  +1.37%  "

# Row 35
$ws.Range("E35").Value = "  +2.32%  "

# Row 37
$ws.Range("D37").Value = "3.675.11"
$ws.Range("E37").Value = "  -1.07%  "

# Row 38
$ws.Range("E38").Value = "  +4.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.63"
$ws.Range("E39").Value = "  +1.70%  "

# Row 40
$ws.Range("E40").Value = "  +7.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.27"
$ws.Range("E41").Value = "  +3.43%  "

# Row 42
$ws.Range("E42").Value = "  +2.60%  "

# Row 43
$ws.Range("E43").Value = "  +1.57%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0423"
$ws.Range("E44").Value = "  +3.85%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  +1.35%  "

# Row 47
$ws.Range("E47").Value = "  +1.36%  "

# Row 48
$ws.Range("E48").Value = "  +5.55%  "

# Row 49
$ws.Range("E49").Value = "  -0.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.59"
$ws.Range("E50").Value = "  +1.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.66"
$ws.Range("E51").Value = "  +4.47%  "
